# Estadisticos Matutinos 15 Oct
# Update the "Rescatables" sheet: reassign the two existing rescatable
# students to a new NC/grade, and add two more rows for the students that
# were previously in rows 2-3 (ANGELES/ROCHA/ADRIAN and FIGUEROA/VAZQUEZ/ORLANDO).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$course = "LECTURA, EXPRESIÓN ORAL Y ESCRITA I"
$group  = "1AV"

# Row 2: new student BAEZ MARCELINO, LUIS EDUARDO
$ws.Cells.Item(2, 1).Value = 21330051920003
$ws.Cells.Item(2, 2).Value = "BAEZ"
$ws.Cells.Item(2, 3).Value = "MARCELINO"
$ws.Cells.Item(2, 4).Value = "LUIS EDUARDO"
$ws.Cells.Item(2, 5).Value = $course
$ws.Cells.Item(2, 6).Value = $group
$ws.Cells.Item(2, 7).Value = 6

# Row 3: new student CARRERA HERNANDEZ, EDGAR FLORENCIO
$ws.Cells.Item(3, 1).Value = 21330051920006
$ws.Cells.Item(3, 2).Value = "CARRERA"
$ws.Cells.Item(3, 3).Value = "HERNANDEZ"
$ws.Cells.Item(3, 4).Value = "EDGAR FLORENCIO"
$ws.Cells.Item(3, 5).Value = $course
$ws.Cells.Item(3, 6).Value = $group
$ws.Cells.Item(3, 7).Value = 6

# Row 4: previously row 2 - ANGELES ROCHA, ADRIAN
$ws.Cells.Item(4, 1).Value = 21330051920002
$ws.Cells.Item(4, 2).Value = "ANGELES"
$ws.Cells.Item(4, 3).Value = "ROCHA"
$ws.Cells.Item(4, 4).Value = "ADRIAN"
$ws.Cells.Item(4, 5).Value = $course
$ws.Cells.Item(4, 6).Value = $group
$ws.Cells.Item(4, 7).Value = 6

# Row 5: previously row 3 - FIGUEROA VAZQUEZ, ORLANDO
$ws.Cells.Item(5, 1).Value = 20330051920291
$ws.Cells.Item(5, 2).Value = "FIGUEROA"
$ws.Cells.Item(5, 3).Value = "VAZQUEZ"
$ws.Cells.Item(5, 4).Value = "ORLANDO"
$ws.Cells.Item(5, 5).Value = $course
$ws.Cells.Item(5, 6).Value = $group
$ws.Cells.Item(5, 7).Value = 6
